$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update selection to F6
[void]$ws.Range("F6").Select()

# Append duplicated company rows (17-46)
$data = @(
    @(17, "El Tule Properties Llc", "803 Brandon Drive", "Seagoville", "TX", "13 White Rock Trl", "Allen", "Tx", 215),
    @(18, "True North Property Owner", "808 Dakota Lane", "Seagoville", "TX", "PO Box 4090", "Scottsdale", "Az", 216),
    @(19, "Sfr 1 2021 1 Borrower Llc", "1625 Emily Lane", "Seagoville", "TX", "1508 Brookhollow Dr", "Santa Ana", "Ca", 217),
    @(20, "Safari One Asset Company Llc", "1611 Emily Lane", "Seagoville", "TX", "5001 Plaza On The Lk Ste 200", "Austin", "Tx", 218),
    @(21, "Home Sfr Borrower Llc", "813 Huddleston Court", "Seagoville", "TX", "3505 Kroger Blvd", "Duluth", "Ga", 219),
    @(22, "Amc Homes Llc", "1002 Fawn Meadow Drive", "Seagoville", "TX", "415 Estate Ln", "Terrell", "Tx", 220),
    @(23, "Rental Transition Llc", "1519 La Fonda Drive", "Seagoville", "TX", "8765 Cleaver Ln", "Terrell", "Tx", 221),
    @(24, "Smfamilyholdings Llc", "1516 La Fonda Drive", "Seagoville", "TX", "7505 Vista Ridge Ct", "Garland", "Tx", 222),
    @(25, "City &amp; Cntry Hms Of America Llc", "1514 La Fonda Drive", "Seagoville", "TX", "5576 N Interstate Highway 45", "Ennis", "Tx", 223),
    @(26, "City &amp; Cntry Hms Of America Llc", "1512 La Fonda Drive", "Seagoville", "TX", "5576 N Interstate Highway 45", "Ennis", "Tx", 224),
    @(27, "City &amp; Cntry Hms Of America Llc", "1510 La Fonda Drive", "Seagoville", "TX", "5576 N Interstate Highway 45", "Ennis", "Tx", 225),
    @(28, "Md Thompson Leasing Company 4", "1508 La Fonda Drive", "Seagoville", "TX", "1918 Seagoville Rd", "Seagoville", "Tx", 226),
    @(29, "City &amp; Cntry Hms Of America Llc", "1506 La Fonda Drive", "Seagoville", "TX", "5576 N Interstate Highway 45", "Ennis", "Tx", 227),
    @(30, "Md Thompson Leasing Company 4", "1504 La Fonda Drive", "Seagoville", "TX", "1918 Seagoville Rd", "Seagoville", "Tx", 228),
    @(31, "City &amp; Cntry Hms Of America Llc", "1502 La Fonda Drive", "Seagoville", "TX", "5576 N Interstate Highway 45", "Ennis", "Tx", 229),
    @(32, "El Tule Properties Llc", "803 Brandon Drive", "Seagoville", "TX", "13 White Rock Trl", "Allen", "Tx", 329),
    @(33, "True North Property Owner", "808 Dakota Lane", "Seagoville", "TX", "PO Box 4090", "Scottsdale", "Az", 330),
    @(34, "Sfr 1 2021 1 Borrower Llc", "1625 Emily Lane", "Seagoville", "TX", "1508 Brookhollow Dr", "Santa Ana", "Ca", 331),
    @(35, "Safari One Asset Company Llc", "1611 Emily Lane", "Seagoville", "TX", "5001 Plaza On The Lk Ste 200", "Austin", "Tx", 332),
    @(36, "Home Sfr Borrower Llc", "813 Huddleston Court", "Seagoville", "TX", "3505 Kroger Blvd", "Duluth", "Ga", 333),
    @(37, "Amc Homes Llc", "1002 Fawn Meadow Drive", "Seagoville", "TX", "415 Estate Ln", "Terrell", "Tx", 334),
    @(38, "Rental Transition Llc", "1519 La Fonda Drive", "Seagoville", "TX", "8765 Cleaver Ln", "Terrell", "Tx", 335),
    @(39, "Smfamilyholdings Llc", "1516 La Fonda Drive", "Seagoville", "TX", "7505 Vista Ridge Ct", "Garland", "Tx", 336),
    @(40, "City &amp; Cntry Hms Of America Llc", "1514 La Fonda Drive", "Seagoville", "TX", "5576 N Interstate Highway 45", "Ennis", "Tx", 337),
    @(41, "City &amp; Cntry Hms Of America Llc", "1512 La Fonda Drive", "Seagoville", "TX", "5576 N Interstate Highway 45", "Ennis", "Tx", 338),
    @(42, "City &amp; Cntry Hms Of America Llc", "1510 La Fonda Drive", "Seagoville", "TX", "5576 N Interstate Highway 45", "Ennis", "Tx", 339),
    @(43, "Md Thompson Leasing Company 4", "1508 La Fonda Drive", "Seagoville", "TX", "1918 Seagoville Rd", "Seagoville", "Tx", 340),
    @(44, "City &amp; Cntry Hms Of America Llc", "1506 La Fonda Drive", "Seagoville", "TX", "5576 N Interstate Highway 45", "Ennis", "Tx", 341),
    @(45, "Md Thompson Leasing Company 4", "1504 La Fonda Drive", "Seagoville", "TX", "1918 Seagoville Rd", "Seagoville", "Tx", 342),
    @(46, "City &amp; Cntry Hms Of America Llc", "1502 La Fonda Drive", "Seagoville", "TX", "5576 N Interstate Highway 45", "Ennis", "Tx", 343)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 12).Value = $row[8]
}
